$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Chapter 6 ("Functions") exercises completed count
$ws.Range("D7").Value = 56

# Re-enter the I column formula as a single fill so Excel stores it as one shared formula group
$ws.Range("I2:I20").Formula = "=H2+G2-1"

# Update the active selection to match the recorded sheet view state
$ws.Range("D8:I8").Select()
